# Insert a new data row at row 69 (above the current row 69), shifting all
# existing rows 69-151 down to 70-152. Then populate the new row 69 with the
# new "Poroto verde" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 69; Excel copies formatting (e.g. the date
# number format on column D) from the surrounding rows automatically.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new record's values.
$ws.Range("A69").Value = 10
$ws.Range("B69").Value = "Vega Modelo de Temuco"
$ws.Range("C69").Value = "La Araucanía"
$ws.Range("D69").Value = 44792
$ws.Range("E69").Value = 9
$ws.Range("F69").Value = 100112031
$ws.Range("G69").Value = "Poroto verde"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 45
$ws.Range("K69").Value = 33000
$ws.Range("L69").Value = 33000
$ws.Range("M69").Value = 33000
$ws.Range("N69").Value = "`$/malla 25 kilos"
$ws.Range("O69").Value = "Provincia de Limarí"
$ws.Range("P69").Value = 1320
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"
